# Automation Script for InvalidLogin TestCase with test data
#
# Adds a new "InvalidLogin" worksheet (after the existing "ValidLogin"
# sheet) containing the same username/password headers plus an invalid
# set of test-data credentials, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current "ValidLogin" sheet so the
# tab order becomes ValidLogin, InvalidLogin.
$validLogin = $wb.Worksheets.Item("ValidLogin")
$invalidLogin = $wb.Worksheets.Add($null, $validLogin)
$invalidLogin.Name = "InvalidLogin"

# Header row - same labels as the ValidLogin sheet.
$invalidLogin.Range("A1").Value = "username"
$invalidLogin.Range("B1").Value = "password"

# Invalid-login test data.
$invalidLogin.Range("A2").Value = "abcd"
$invalidLogin.Range("B2").Value = "xyz"

# Leave the new sheet active with B2 selected, mirroring the authored file.
[void]$invalidLogin.Range("B2").Select()
